# "added address column in excel sheets"
# Insert a new column F ("Address") immediately before the existing
# "District" column, which shifts right from F to G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting before column F pushes the old F column ("District") to G,
# leaving a blank column F ready for the new "Address" data.
$ws.Columns("F:F").Insert()

# New column F header (row 2) and per-row address values (rows 4-44).
# Row 3 (the "NAME"/sub-header row) has no address and stays blank.
$addressByRow = [ordered]@{
    2 = 'Address'
    4 = 'N S F H S KulagodMudalagiChikkodi'
    5 = 'G H S AvaragolHukkeriEPU'
    6 = 'New High School KempwadAthani'
    7 = 'S D High SchoolSankeshwarHukkeri'
    8 = 'S A P V M KanagalaHukeri'
    9 = 'S B S Kanya ShalaNipani'
    10 = 'S S High School KhanagaonGokak'
    11 = 'Govt. High School NeermanviManvi'
    12 = 'G G H S KavitalManvi'
    13 = 'Govt. High SchoolTalakatanalMudalagiGokak'
    14 = 'G H S KhanagaonGokak'
    15 = 'Shri ParshwamatiKanya VidyalayAkol'
    16 = 'Govt. P U College (High School Section) MajalattiChikkodi'
    17 = 'Govt. High SchoolNeeralakeriLingasugur'
    18 = 'G G H S ShindikurbetGokak'
    19 = 'G H S Nadi – IngalagaonAthani'
    20 = 'K S S High School JugulAthani'
    21 = 'Bhirdi High School BhirdiBhirdiRaibag'
    22 = 'Govt. P U College GokakGokak'
    23 = 'G H S BasapurHukkeri'
    24 = 'Govt. High School Matamari'
    25 = 'Govt. High SchoolK E B Colony'
    26 = 'Govt. High School HalalliAthani'
    27 = 'Govt. High School MavinabhaviLingasugar'
    28 = 'G P U College (High School Section) NaganurGokakMudalagi'
    29 = 'Govt. High School SaavasuddiRaibagChikkodi'
    30 = 'G K H P SchoolWadagolNippaniChikkodi'
    31 = 'G H S GurlapurMudalagi'
    32 = 'Govt. High School HampanalSindanuru'
    33 = 'Govt. High SchoolNagaramunnoliChikkodi'
    34 = 'Ajitkumar Baane High School HandigundRaibaag'
    35 = 'G H P S PudakalakattiGokak'
    36 = 'G H S Udbal (U)Sindhanur'
    37 = 'G H SchoolNallanatti Mudalgi'
    38 = 'T G TK H P SPattanakudiChikkodi'
    39 = 'U G H P S JanekalManvi'
    40 = 'Govt High School SangapurManvi'
    41 = 'G P U H S UrubinahattiGokak'
    42 = 'Govt. High SchoolKerurChikkodi'
    43 = 'G P U (High School Section) Hutti Gold Mines'
    44 = 'S M S ShindikurbetGokak'
}

foreach ($row in $addressByRow.Keys) {
    $ws.Range("F$row").Value = $addressByRow[$row]
}

# Row 3 (between the "Address" header and the first data row) has no
# address text in the source data, so its new F cell is left blank.
$ws.Range("F3").Value = ""

